# Apply cryptos list update (generated from OOXML diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a pure numeric-looking string that must stay
# textual (keeps trailing zeros / avoids scientific notation), matching
# the original inlineStr representation.
$textCells = @("D12", "D13", "D29", "D32", "D40")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "73.215.11"
$ws.Range("E2").Value = "  +3.20%  "
$ws.Range("D3").Value = "3.994.73"
$ws.Range("E3").Value = "  +1.35%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "597.25"
$ws.Range("E5").Value = "  +11.58%  "
$ws.Range("D6").Value = "163.77"
$ws.Range("E6").Value = "  +11.03%  "
$ws.Range("D7").Value = "0.685"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D9").Value = "0.751"
$ws.Range("E9").Value = "  +1.95%  "
$ws.Range("D10").Value = "0.169"
$ws.Range("E10").Value = "  +2.59%  "
$ws.Range("D11").Value = "54.63"
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("D12").Value = "0.0000321"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("D13").Value = "11.00"
$ws.Range("E13").Value = "  +3.95%  "
$ws.Range("D14").Value = "4.619.55"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").Value = "3.990.44"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("E16").Value = "  +9.34%  "
$ws.Range("D17").Value = "14.16"
$ws.Range("E17").Value = "  +2.28%  "
$ws.Range("D18").Value = "20.42"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").Value = "72.771.79"
$ws.Range("E20").Value = "  +2.73%  "
$ws.Range("D21").Value = "439.03"
$ws.Range("E21").Value = "  +4.34%  "
$ws.Range("D22").Value = "4.75"
$ws.Range("E22").Value = "  +12.71%  "
$ws.Range("D23").Value = "96.47"
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("D24").Value = "3.45"
$ws.Range("E24").Value = "  -4.34%  "
$ws.Range("D25").Value = "14.34"
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("D26").Value = "4.36"
$ws.Range("E26").Value = "  +14.72%  "
$ws.Range("D27").Value = "11.32"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("D29").Value = "10.40"
$ws.Range("E29").Value = "  -2.52%  "
$ws.Range("D30").Value = "36.33"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "7.82"
$ws.Range("E31").Value = "  -0.81%  "
$ws.Range("D32").Value = "13.80"
$ws.Range("E32").Value = "  +3.49%  "
$ws.Range("D33").Value = "0.131"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").Value = "48.33"
$ws.Range("E34").Value = "  -4.97%  "
$ws.Range("D35").Value = "672.03"
$ws.Range("E35").Value = "  -1.78%  "
$ws.Range("D36").Value = "70.98"
$ws.Range("E36").Value = "  +8.60%  "
$ws.Range("D37").Value = "0.0₃0904"
$ws.Range("E37").Value = "  +11.36%  "
$ws.Range("D38").Value = "0.437"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("B39").Value = "ThetaToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D39").Value = "3.36"
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.145"
$ws.Range("E41").Value = "  -1.65%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "3.34"
$ws.Range("E42").Value = "  +5.24%  "
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "0.0491"
$ws.Range("E44").Value = "  +2.24%  "
$ws.Range("D45").Value = "10.67"
$ws.Range("E45").Value = "  +6.93%  "
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("D48").Value = "2.63"
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("D49").Value = "2.895.26"
$ws.Range("E49").Value = "  +10.19%  "
$ws.Range("D50").Value = "3.06"
$ws.Range("E50").Value = "  +2.22%  "
$ws.Range("D51").Value = "3.41"
$ws.Range("E51").Value = "  +4.85%  "

# Restore default (General) styling on the forced-text cells so the
# stored cell style matches the original (unstyled) cells.
foreach ($c in $textCells) {
    $ws.Range($c).Style = "Normal"
}
